$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.491.91"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "2.422.97"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.60%  "
$ws.Range("E7").Value = "  +2.23%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.23%  "
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "2.800.59"
$ws.Range("D16").Value = "2.462.57"
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("E17").Value = "  +5.28%  "
$ws.Range("D18").Value = "44.351.75"
$ws.Range("E18").Value = "  +3.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").Value = "0.0₃0919"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.75%  "
$ws.Range("E24").Value = "  +5.57%  "
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("E32").Value = "  +17.38%  "
$ws.Range("E33").Value = "  +10.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0770"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.75%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.26%  "
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0289"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("D45").Value = "1.943.68"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.95%  "
